$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.268.82'
$ws.Range("D3").Value = '3.746.62'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("D7").Value = '3.743.73'
$ws.Range("E7").Value = '  +0.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("E10").Value = '  +5.41%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.25'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("E14").Value = '  +2.27%  '
$ws.Range("D15").Value = '4.372.41'
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("D16").Value = '3.746.80'
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("D17").Value = '69.227.57'
$ws.Range("E17").Value = '  +0.88%  '
$ws.Range("E18").Value = '  +2.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.81%  '
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +12.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '493.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.729'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000150'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.60%  '
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.66%  '
$ws.Range("E31").Value = '  +3.52%  '
$ws.Range("E32").Value = '  +0.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.66'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").Value = '3.892.83'
$ws.Range("E34").Value = '  +0.56%  '
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("D36").Value = '3.678.83'
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.99'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.14%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.140'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("E41").Value = '  +0.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '423.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.46%  '
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").Value = '2.796.20'
$ws.Range("E50").Value = '  +1.81%  '
$ws.Range("E51").Value = '  +0.64%  '
